$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the old "selectors" lesson titles from G12:G15 - the lessons are
# being replaced by the new "CSS Text" lesson entries appended below.
$ws.Range("G12:G15").ClearContents()

# The F6:F18 merged block (topic banner) is no longer centered/merged now
# that the list below it has grown - unmerge it and drop the centered
# alignment back to general.
$ws.Range("F6:F18").UnMerge()
$ws.Range("F6:F18").HorizontalAlignment = 1

# Append the two new "CSS Text" lessons at the bottom of the list.
$ws.Range("F71").Value = 4
$ws.Range("G71").Value = "white-space?"
$ws.Range("F72").Value = 5
$ws.Range("G72").Value = "list-style"

# Scroll the view down to the newly added rows and select the last entry,
# matching where the author was working when the sheet was saved.
$ws.Range("F70").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 43
